$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 6-9 (Q4, Q5, Q6, Q7) entirely - they are no longer part of the table
$ws.Range("A6:G9").EntireRow.Delete()

# Update the remaining data values (rows 2-5) with the corrected/re-simulated figures
$ws.Range("B2").Value = 0.4532102921767785
$ws.Range("C2").Value = 0.9560173308380451
$ws.Range("D2").Value = 2.911849634940581
$ws.Range("E2").Value = 1.706414262405404
$ws.Range("F2").Value = 1.707231247553019
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.5397191998186066
$ws.Range("C3").Value = 1.131880666064611
$ws.Range("D3").Value = 3.375819745879244
$ws.Range("E3").Value = 1.837340400110781
$ws.Range("F3").Value = 1.851282717837196
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = 0.01424612316995905
$ws.Range("C4").Value = 0.3434453827890073
$ws.Range("D4").Value = 0.1563237768149852
$ws.Range("E4").Value = 0.395378017617299
$ws.Range("F4").Value = 0.432833674461142
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 0.1208118526715265
$ws.Range("C5").Value = 0.1596663445948482
$ws.Range("D5").Value = 0.04008884534220743
$ws.Range("E5").Value = 0.2002219901564447
$ws.Range("F5").Value = 0.2258023099805705
$ws.Range("G5").Value = 2
